# Insert a new data row right before the existing row 399, which shifts the
# existing rows 399-423 down to 400-424 (matches the diff: dimension goes
# from A1:T423 to A1:T424, and every row from 399 onward now holds the data
# that used to belong to the row above it).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(399).Insert()

# Populate the newly inserted row 399 with its new values.
$ws.Cells.Item(399, 1).Value = 9
$ws.Cells.Item(399, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(399, 3).Value = "Metropolitana"
$ws.Cells.Item(399, 4).Value = 44706
$ws.Cells.Item(399, 5).Value = 13
$ws.Cells.Item(399, 6).Value = "Fruta"
$ws.Cells.Item(399, 7).Value = 100108
$ws.Cells.Item(399, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(399, 9).Value = 100108002
$ws.Cells.Item(399, 10).Value = "Mango"
$ws.Cells.Item(399, 11).Value = "Sin especificar"
$ws.Cells.Item(399, 12).Value = "Primera"
$ws.Cells.Item(399, 13).Value = 570
$ws.Cells.Item(399, 14).Value = 7000
$ws.Cells.Item(399, 15).Value = 7500
$ws.Cells.Item(399, 16).Value = 7246
$ws.Cells.Item(399, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(399, 18).Value = "Brasil"
$ws.Cells.Item(399, 19).Value = 1812
$ws.Cells.Item(399, 20).Value = 4
